$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Administración delegada a través de Corporaciones Municipales. Éstas son unidades responsables de la administración global de recursos humanos, físicos y financieros asociados a los estudiantes que asisten a los establecimientos municipales de la comuna y también están encargadas de asesorar a los Alcaldes y Concejos Municipales en lo relativo a la formulación de políticas en el ámbito educativo.'
$ws.Range("C3").Value = 'Establecimientos educacionales cuya administración ha sido traspasada a las municipalidades del país.'
$ws.Range("C4").Value = 'Establecimientos que pertenecen a sostenedores privados y que reciben una subvención del Estado por cada alumno atendido.'
$ws.Range("C5").Value = 'Establecimientos que pertenecen a particulares, que no reciben subvención del Estado y que son pagados por los padres y apoderados.'
$ws.Range("C6").Value = 'Un sostenedor privado es el que administra o gestiona el establecimiento educacional, que es de propiedad del estado, y que recibe recursos del estado para su funcionamiento. Esta modalidad está restringida tanto a sostenedores de derecho privado sin fines de lucro, como aquellos que imparten educación técnico profesional.'
$ws.Range("C9").Value = 'Es uno de los componentes del Sistema de Aseguramiento de la Calidad de la Educación, cuyos objetivos son evaluar el aprendizaje de los estudiantes y el logro de los otros indicadores de calidad educativa, informar a la comunidad escolar sobre la evaluación de los establecimientos e identificar las necesidades de apoyo, en especial a los establecimientos ordenados en categoría de desempeño Medio- Bajo e Insuficiente. Su construcción se basa en un Índice de Resultados inicial que considera: la distribución de los estudiantes en los Niveles de Aprendizaje (67%), los Indicadores de Desarrollo Personal y Social, los resultados de las pruebas Simce y su progreso en las últimas tres o dos mediciones según corresponda para cada nivel (33%). Luego, este Índice de Resultados se ajusta según las Características de los Estudiantes del establecimiento educacional, por ejemplo, su vulnerabilidad. Finalmente, en base a este nuevo Índice de Resultados final, se clasifica a los establecimientos en Alto, Medio, Medio bajo e Insuficiente.'
$ws.Range("C10").Value = 'Agrupa a establecimientos cuyos estudiantes obtienen resultados que sobresalen respecto de lo esperado, considerando siempre el contexto social de los estudiantes del establecimiento.'
$ws.Range("C11").Value = 'Agrupa a establecimientos cuyos estudiantes obtienen resultados similares a lo esperado, considerando siempre el contexto social de los estudiantes del establecimiento.'
$ws.Range("C12").Value = 'Agrupa establecimientos cuyos estudiantes obtienen resultados por debajo de lo esperado, considerando siempre el contexto social de los estudiantes del establecimiento.'
$ws.Range("C13").Value = 'Agrupa a establecimientos cuyos estudiantes obtienen resultados muy por debajo de lo esperado, considerando siempre el contexto social de los estudiantes del establecimiento.'
$ws.Range("C14").Value = 'Los Estándares de Aprendizaje describen lo que los estudiantes deben saber y poder hacer para demostrar si alcanzan los objetivos de aprendizaje estipulados en el currículo vigente.'
$ws.Range("C15").Value = 'Los estudiantes que quedan clasificados en este nivel no logran demostrar consistentemente que han adquirido los conocimientos y habilidades más elementales estipulados en el currículo para el período evaluado.'
$ws.Range("C16").Value = 'Los estudiantes que alcanzan este Nivel de Aprendizaje han logrado lo exigido en el currículo de manera parcial. Esto implica demostrar que han adquirido los conocimientos y habilidades más elementales estipulados para el período evaluado.'
$ws.Range("C17").Value = 'Los estudiantes que alcanzan este Nivel de Aprendizaje han logrado lo exigido en el currículo de manera satisfactoria. Esto implica demostrar que han adquirido los conocimientos y habilidades básicos estipulados para el período evaluado.'
$ws.Range("C18").Value = 'Un estudiante que se siente capaz académicamente y que está motivado por el estudio, es más probable que se interese e invierta esfuerzo en las actividades escolares.                       Este indicador considera la percepción y valoración de los estudiantes en relación con su capacidad de aprender y por otra parte las percepciones y actitudes que tienen los estudiantes hacia el aprendizaje y el logro académico.'
$ws.Range("C19").Value = 'El clima de convivencia escolar afecta el bienestar y desarrollo socioafectivo de los estudiantes e impacta en su conducta, disposición y rendimiento durante las actividades escolares. Considera las percepciones y las actitudes que tienen los estudiantes, docentes y padres y apoderados con respecto a la presencia de un ambiente de respeto, organizado y seguro.'
$ws.Range("C20").Value = 'Un clima participativo y con un mayor sentido de pertenencia a la escuela mejora el compromiso y las disposición de los estudiantes, docentes, padres y apoderados hacia la mejora de los aprendizajes. Este indicador considera las percepciones y las actitudes que declaran los docentes, estudiantes, padres y apoderados en los cuestionarios que se aplican durante las pruebas Simce.'
$ws.Range("C21").Value = 'La práctica de actividad física y una alimentación balanceada mejora la capacidad de aprendizaje de los estudiantes, ya que incrementa la capacidad para resolver tareas difíciles, la concentración y la memoria. También, disminuye la eventual ansiedad y estrés ante las evaluaciones. Este indicador considera las actitudes y conductas declaradas de los estudiantes en relación con la vida saludable, también sus percepciones sobre el grado en que el establecimiento promueve hábitos de alimentación sana, de vida activa y de autocuidado.'

$ws.Range("C23").Select() | Out-Null
